$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Shift the two existing data rows up (row3->row2, row4->row3) and
#     rebuild row 1 as a proper header row, then add a new blank row 4. ---

# Row 2 (was row 3): Kubel
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 107500
$ws.Range("C2").Value = "Kubel"
$ws.Range("D2").Value = 1900
$ws.Range("E2").Value = 1976
$ws.Range("F2").Value = 18
$ws.Range("G2").Value = 1.48
$ws.Range("H2").Value = 1.37
$ws.Range("I2").Value = 1.4
$ws.Range("J2").Value = 1.66
$ws.Range("K2").Value = 3.06

# Row 3 (was row 4): Wasserauen
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 107400
$ws.Range("C3").Value = "Wasserauen"
$ws.Range("D3").Value = 1905
$ws.Range("E3").Value = 2005
$ws.Range("F3").Value = 1.3
$ws.Range("G3").Value = 2.5
$ws.Range("H3").Value = 2.5
$ws.Range("I3").Value = 1.48
$ws.Range("J3").Value = 6.22
$ws.Range("K3").Value = 7.7

# New header row 1
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# Row 4 becomes blank (used to be the second data row before the shift)
$ws.Range("A4").Value = ""
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = ""

# Drop the now-superfluous trailing blank row (sheet used to run to row 31).
$ws.Rows("31:31").Delete()

# Update the visible selection to match the authored state (row 2 selected).
$ws.Range("A2:K2").Select()

Write-Output "done"
